$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,8
$data[0,0] = [double]0
$data[0,1] = "walkingToRunning"
$data[0,2] = [double]-7.487135410308838
$data[0,3] = [double]2.468842267990112
$data[0,4] = [double]-10.7623405456543
$data[0,5] = [double]0.00006657902849829999
$data[0,6] = [double]0.5838314890861511
$data[0,7] = [double]0.3130545914173126

$data[1,0] = [double]100
$data[1,1] = "walkingToRunning"
$data[1,2] = [double]5.30696439743042
$data[1,3] = [double]-7.900328636169434
$data[1,4] = [double]-8.124782562255859
$data[1,5] = [double]1.05201518535614
$data[1,6] = [double]1.37851881980896
$data[1,7] = [double]0.5420864224433899

$data[2,0] = [double]200
$data[2,1] = "walkingToRunning"
$data[2,2] = [double]0.1120047569274902
$data[2,3] = [double]-3.523766279220581
$data[2,4] = [double]-6.589697360992432
$data[2,5] = [double]-1.053479909896851
$data[2,6] = [double]-0.3754391372203827
$data[2,7] = [double]0.5234442949295044

$data[3,0] = [double]300
$data[3,1] = "walkingToRunning"
$data[3,2] = [double]-0.1642012596130371
$data[3,3] = [double]-2.677549123764038
$data[3,4] = [double]-5.916802883148193
$data[3,5] = [double]-1.651625990867615
$data[3,6] = [double]-0.3386875092983246
$data[3,7] = [double]2.086187362670898

$data[4,0] = [double]400
$data[4,1] = "walkingToRunning"
$data[4,2] = [double]2.91043758392334
$data[4,3] = [double]-5.157215595245361
$data[4,4] = [double]-1.494926452636719
$data[4,5] = [double]2.480534791946411
$data[4,6] = [double]0.6072673201560974
$data[4,7] = [double]-0.0555269084870815

$data[5,0] = [double]500
$data[5,1] = "walkingToRunning"
$data[5,2] = [double]-4.463288307189941
$data[5,3] = [double]-6.138980388641357
$data[5,4] = [double]-15.2139720916748
$data[5,5] = [double]0.3121890723705292
$data[5,6] = [double]-3.295861721038818
$data[5,7] = [double]-2.403902292251587

$data[6,0] = [double]600
$data[6,1] = "walkingToRunning"
$data[6,2] = [double]-2.665878295898437
$data[6,3] = [double]-4.94788122177124
$data[6,4] = [double]-14.101637840271
$data[6,5] = [double]3.045657634735107
$data[6,6] = [double]-1.371461391448975
$data[6,7] = [double]-3.965579986572266

$data[7,0] = [double]700
$data[7,1] = "walkingToRunning"
$data[7,2] = [double]2.109983444213867
$data[7,3] = [double]-22.90362548828125
$data[7,4] = [double]-9.720759391784668
$data[7,5] = [double]1.784917116165161
$data[7,6] = [double]3.971904993057251
$data[7,7] = [double]-6.264420986175537

$data[8,0] = [double]800
$data[8,1] = "walkingToRunning"
$data[8,2] = [double]-23.78427696228028
$data[8,3] = [double]-3.573103904724121
$data[8,4] = [double]-5.555922031402588
$data[8,5] = [double]-2.070807456970215
$data[8,6] = [double]-2.531001806259156
$data[8,7] = [double]5.346961975097656

$data[9,0] = [double]900
$data[9,1] = "walkingToRunning"
$data[9,2] = [double]-0.0823409557342529
$data[9,3] = [double]6.524863243103027
$data[9,4] = [double]0.961866855621338
$data[9,5] = [double]-5.287373542785645
$data[9,6] = [double]-1.623929142951965
$data[9,7] = [double]4.378103733062744

$data[10,0] = [double]1000
$data[10,1] = "walkingToRunning"
$data[10,2] = [double]-3.332920551300049
$data[10,3] = [double]-42.73205184936523
$data[10,4] = [double]-5.979649543762207
$data[10,5] = [double]-9.357216835021973
$data[10,6] = [double]4.914664268493652
$data[10,7] = [double]2.073936700820923

$data[11,0] = [double]1100
$data[11,1] = "walkingToRunning"
$data[11,2] = [double]7.95712423324585
$data[11,3] = [double]9.703181266784668
$data[11,4] = [double]-9.591499328613279
$data[11,5] = [double]1.831788778305054
$data[11,6] = [double]-9.64909839630127
$data[11,7] = [double]-1.40841281414032

$data[12,0] = [double]1200
$data[12,1] = "walkingToRunning"
$data[12,2] = [double]-12.04059410095215
$data[12,3] = [double]-11.62393760681152
$data[12,4] = [double]-15.216552734375
$data[12,5] = [double]3.237405300140381
$data[12,6] = [double]-1.015663027763367
$data[12,7] = [double]-6.009289741516113

$data[13,0] = [double]1300
$data[13,1] = "walkingToRunning"
$data[13,2] = [double]-30.35504722595215
$data[13,3] = [double]-42.07619094848633
$data[13,4] = [double]14.66377067565918
$data[13,5] = [double]5.363673210144043
$data[13,6] = [double]3.099453449249268
$data[13,7] = [double]-2.63826060295105

$data[14,0] = [double]1400
$data[14,1] = "walkingToRunning"
$data[14,2] = [double]4.246713161468506
$data[14,3] = [double]3.201902866363525
$data[14,4] = [double]-3.634400367736816
$data[14,5] = [double]-5.723599433898926
$data[14,6] = [double]2.653107643127441
$data[14,7] = [double]2.690192222595215

$data[15,0] = [double]1500
$data[15,1] = "walkingToRunning"
$data[15,2] = [double]15.14187526702881
$data[15,3] = [double]-2.368836879730225
$data[15,4] = [double]2.886821746826172
$data[15,5] = [double]-7.246927738189697
$data[15,6] = [double]-1.130711674690247
$data[15,7] = [double]3.536544799804688

$data[16,0] = [double]1600
$data[16,1] = "walkingToRunning"
$data[16,2] = [double]-66.83038330078125
$data[16,3] = [double]-23.19043350219727
$data[16,4] = [double]-35.43113708496094
$data[16,5] = [double]0.5896904468536377
$data[16,6] = [double]12.76406574249268
$data[16,7] = [double]2.221475839614868

$data[17,0] = [double]1700
$data[17,1] = "walkingToRunning"
$data[17,2] = [double]21.56573677062988
$data[17,3] = [double]-0.212137758731842
$data[17,4] = [double]-20.44221496582031
$data[17,5] = [double]-1.113667368888855
$data[17,6] = [double]-3.135539293289185
$data[17,7] = [double]2.971954584121704

$data[18,0] = [double]1800
$data[18,1] = "walkingToRunning"
$data[18,2] = [double]43.28559112548828
$data[18,3] = [double]-13.05164337158203
$data[18,4] = [double]-13.05863285064697
$data[18,5] = [double]6.756506443023682
$data[18,6] = [double]-7.358247756958008
$data[18,7] = [double]-5.698232650756836

$data[19,0] = [double]1900
$data[19,1] = "walkingToRunning"
$data[19,2] = [double]30.09039688110352
$data[19,3] = [double]-34.86870574951172
$data[19,4] = [double]6.082807540893555
$data[19,5] = [double]4.204665184020996
$data[19,6] = [double]-1.005010485649109
$data[19,7] = [double]-2.928545236587524

$data[20,0] = [double]2000
$data[20,1] = "walkingToRunning"
$data[20,2] = [double]-2.691319465637207
$data[20,3] = [double]10.76634311676025
$data[20,4] = [double]-5.473065376281738
$data[20,5] = [double]-1.400756239891052
$data[20,6] = [double]0.9566740393638612
$data[20,7] = [double]2.680604934692383

$data[21,0] = [double]2100
$data[21,1] = "walkingToRunning"
$data[21,2] = [double]11.3234281539917
$data[21,3] = [double]-40.29625701904297
$data[21,4] = [double]38.12360763549805
$data[21,5] = [double]-11.39346981048584
$data[21,6] = [double]-1.724596619606018
$data[21,7] = [double]13.37479496002197

$data[22,0] = [double]2200
$data[22,1] = "walkingToRunning"
$data[22,2] = [double]-7.852145671844482
$data[22,3] = [double]6.264562606811523
$data[22,4] = [double]-14.15214920043945
$data[22,5] = [double]8.482767105102539
$data[22,6] = [double]-0.4398876428604126
$data[22,7] = [double]-8.189353942871094

$data[23,0] = [double]2300
$data[23,1] = "walkingToRunning"
$data[23,2] = [double]-4.107211589813232
$data[23,3] = [double]-2.712513208389282
$data[23,4] = [double]-20.56048202514648
$data[23,5] = [double]-0.3088601231575012
$data[23,6] = [double]-7.590475082397461
$data[23,7] = [double]2.290185451507568

$data[24,0] = [double]2400
$data[24,1] = "walkingToRunning"
$data[24,2] = [double]-8.579601287841797
$data[24,3] = [double]-15.95286655426025
$data[24,4] = [double]-10.86754608154297
$data[24,5] = [double]6.10509729385376
$data[24,6] = [double]-0.359460175037384
$data[24,7] = [double]-4.315519332885742

$data[25,0] = [double]2500
$data[25,1] = "walkingToRunning"
$data[25,2] = [double]0.7219026684761047
$data[25,3] = [double]2.51579213142395
$data[25,4] = [double]17.58181953430176
$data[25,5] = [double]-1.245227575302124
$data[25,6] = [double]-1.861483097076416
$data[25,7] = [double]1.978062987327576

$data[26,0] = [double]2600
$data[26,1] = "walkingToRunning"
$data[26,2] = [double]-6.065989017486572
$data[26,3] = [double]16.71288681030273
$data[26,4] = [double]-1.314104557037354
$data[26,5] = [double]-1.619668006896973
$data[26,6] = [double]1.064798355102539
$data[26,7] = [double]1.856090188026428

$data[27,0] = [double]2700
$data[27,1] = "walkingToRunning"
$data[27,2] = [double]-11.12415027618408
$data[27,3] = [double]-78.97219848632812
$data[27,4] = [double]36.16990280151367
$data[27,5] = [double]-10.51302814483643
$data[27,6] = [double]-17.02712059020996
$data[27,7] = [double]-4.28569221496582

$data[28,0] = [double]2800
$data[28,1] = "walkingToRunning"
$data[28,2] = [double]-0.5275765657424927
$data[28,3] = [double]10.35773277282715
$data[28,4] = [double]-23.80691909790039
$data[28,5] = [double]9.720071792602541
$data[28,6] = [double]-2.732869386672974
$data[28,7] = [double]-1.917076587677002

$data[29,0] = [double]2900
$data[29,1] = "walkingToRunning"
$data[29,2] = [double]-18.0826530456543
$data[29,3] = [double]-1.298346519470215
$data[29,4] = [double]-14.97142791748047
$data[29,5] = [double]0.3920839130878448
$data[29,6] = [double]-4.407997608184815
$data[29,7] = [double]1.71547520160675

$ws.Range("A2:H31").Value2 = $data